# Jenkins parameterized build completed
#
# Adds a "Runmode" column (E) to the addCustomerTest sheet, populates it
# with alternating y/n values for the existing data rows, and makes
# addCustomerTest the active sheet/selection (it was test_suite before).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("addCustomerTest")

# New "Runmode" column with per-row run flags.
$ws1.Range("E1").Value = "Runmode"
$ws1.Range("E2").Value = "y"
$ws1.Range("E3").Value = "n"
$ws1.Range("E4").Value = "y"
$ws1.Range("E5").Value = "n"

# Make addCustomerTest the active sheet and select the last written cell,
# matching the new workbook/sheet view state.
$ws1.Activate()
$ws1.Range("E5").Select()
